$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Designator D72) previously had an empty DESCRIPTION (column C).
# Fill in the diode order-code / description, matching the formatting
# already used by the surrounding data rows (wrap text + border, same as
# the other filled cells in column C).
$ws.Range("C7").Value = "WE-TVS TVS Diode, High Speed Series, SOT23-6L, VDC = 5V"
$ws.Range("C8").Copy()
$ws.Range("C7").PasteSpecial(-4122)

# The longer description text now wraps onto two lines, so the row grows
# taller (matches row 8/9's wrapped-text height).
$ws.Rows.Item(7).RowHeight = 28.8

# Print scale was nudged down slightly.
$ws.PageSetup.Zoom = 84
